$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "63.860.12"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  +0.43%  "
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "2.635.82"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  +0.98%  "
$ws.Range("E4").Value = "  -0.10%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "578.65"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +0.92%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "156.90"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +1.62%  "
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.631"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  -0.44%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -1.96%  "
$ws.Range("E10").Value = "  +1.32%  "
$ws.Range("E11").Value = "  +0.02%  "
$ws.Range("E12").Value = "  +1.04%  "
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "28.71"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  +1.73%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "3.112.64"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  +0.91%  "
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "0.0000184"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  +0.98%  "
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "63.776.36"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  +0.52%  "
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "2.631.89"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  +0.81%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "7.74"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  +3.47%  "
$ws.Range("E20").Value = "  -1.75%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "344.48"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  +0.59%  "
$ws.Range("E22").Value = "  +0.17%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "68.35"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  +2.03%  "
$ws.Range("E24").Value = "  +8.47%  "
$ws.Range("E25").Value = "  +4.89%  "
$ws.Range("B26").Value = "InternetComputer(DFINITY)"
$ws.Range("C26").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "9.28"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  +0.36%  "
$ws.Range("B27").Value = "Fetch.AI"
$ws.Range("C27").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "1.64"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  +5.15%  "
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "580.01"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  +2.38%  "
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "8.23"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  +5.19%  "
$ws.Range("E30").Value = "  +0.84%  "
$ws.Range("E31").Value = "  -0.11%  "
$ws.Range("E32").Value = "  +0.23%  "
$ws.Range("E33").Value = "  +2.47%  "
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "6.62"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  +3.08%  "
$ws.Range("E35").Value = "  +3.38%  "
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.402"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  -1.13%  "
$ws.Range("E37").Value = "  -0.13%  "
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  -0.17%  "
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "1.91"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  +2.99%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "152.76"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  +0.87%  "
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "2.57"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  +9.42%  "
$ws.Range("E42").Value = "  +0.00%  "
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "162.72"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  +4.61%  "
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "24.28"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  +6.47%  "
$ws.Range("E45").Value = "  -0.46%  "
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "0.0588"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -0.53%  "
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "0.634"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  +1.18%  "
$ws.Range("E48").Value = "  -1.16%  "
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "0.0248"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  -0.47%  "
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "19.06"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  +0.65%  "
$ws.Range("E51").Value = "  +1.84%  "
